$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 2) below the existing header row.
# Column C holds a text value that looks numeric ("2"), so it is entered
# with a leading apostrophe to force text storage instead of a number.
$ws.Range("A2").Value = "Erinna D. Brodsky"
$ws.Range("B2").Value = "Medium"
$ws.Range("C2").Value = "'2"
$ws.Range("D2").Value = 3.77
$ws.Range("E2").Value = 3.96
$ws.Range("F2").Value = "tres"

# Add a totals row (row 3). A3, B3 and F3 are blank text cells (present but
# empty), so a bare apostrophe is used to create an empty-string text cell
# instead of leaving the cell completely untouched.
$ws.Range("A3").Value = "'"
$ws.Range("B3").Value = "'"
$ws.Range("C3").Value = "Total"
$ws.Range("D3").Value = 3.77
$ws.Range("E3").Value = 3.96
$ws.Range("F3").Value = "'"
